$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("articels")
$ws.Activate()

# Clear cells that are removed in the diff
$ws.Range("A3").Value = ""
$ws.Range("A6").Value = ""

# Column A updates / additions (rows 9-19)
$ws.Range("A9").Value  = "/add_vendor"
$ws.Range("A10").Value = "/addven"
$ws.Range("A11").Value = "/addven"
$ws.Range("A12").Value = "/addven"
$ws.Range("A13").Value = "/addven"
$ws.Range("A14").Value = "/addven"
$ws.Range("A15").Value = "/addven"
$ws.Range("A16").Value = "/addven"
$ws.Range("A17").Value = "gfdgfdgfd"
$ws.Range("A18").Value = "жив?"
$ws.Range("A19").Value = "lol"

# Column B additions (rows 2-11)
$ws.Range("B2").Value  = "fsdfsd"
$ws.Range("B3").Value  = "kjk"
$ws.Range("B4").Value  = "realy&"
$ws.Range("B5").Value  = "lol"
$ws.Range("B6").Value  = "gfdgfdg"
$ws.Range("B7").Value  = "добавь"
$ws.Range("B8").Value  = "fdsf"
$ws.Range("B9").Value  = "u"
$ws.Range("B10").Value = "f"
$ws.Range("B11").Value = "fd"

$ws.Range("A7").Select()
